$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Remove the stray _GoBack bookmark from the title line.
#    (A new _GoBack bookmark will be added later, at the end of the
#    "Total Outlier:" paragraph, mirroring the diff.)
# ---------------------------------------------------------------
$gb = $d.Bookmarks("_GoBack")
$gb.Delete()

# ---------------------------------------------------------------
# 2. "95% confidence Interval (Min/Max)  ~1231 / ~3333  Tab is < ~p0005 V"
#    becomes
#    "95% confidence Interval <tab>{{ c_1 }} / {{ c_2 }}  Tab is < {{ c_3 }} V"
# ---------------------------------------------------------------

# 2a. Trim "confidence Interval (Min/Max)" down to "confidence Interval ",
#     then add a fresh tab (bold, sz 18) right after it -- matching the
#     paragraph's "95% " run formatting.
$rng = $d.Content
$rng.Find.Execute("confidence Interval (Min/Max)", $false, $false, $false, $false, $false, $true, 1, $false, "confidence Interval ", 2) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter([char]9)
$rng.Font.Bold = $true
$rng.Font.Size = 9

# 2b. Right after the existing (pre-tab) tab run, insert the "{{ c_1 }}"
#     placeholder using the un-bolded, SimSun-east-asia run formatting
#     that the rest of the line uses.
$rng = $d.Content
$rng.Find.Execute("~1231 / ~3333") | Out-Null
$rng.Collapse(1)
$rng.InsertBefore("{{ c_1 }}")
$rng.Font.Bold = $false
$rng.Font.Size = 10
$rng.Font.NameFarEast = "SimSun"

# 2c. Replace the old "~1231 / ~3333" text with " / {{ c_2 }}" and drop
#     the bold attribute from that run.
$rng = $d.Content
$rng.Find.Execute("~1231 / ~3333", $false, $false, $false, $false, $false, $true, 1, $false, " / {{ c_2 }}", 2) | Out-Null
$rng.Font.Bold = $false

# 2d. "~p0005" -> "{{ c_3 }}" (keep its existing, non-bold formatting).
$rng = $d.Content
$rng.Find.Execute("~p0005", $false, $false, $false, $false, $false, $true, 1, $false, "{{ c_3 }}", 2) | Out-Null

# ---------------------------------------------------------------
# 3. Insert a brand-new "Outlier Min/Max" paragraph right after the
#    confidence-interval paragraph (and before "Total Samples Outside
#    Range:"). It inherits the bordered-box paragraph formatting of its
#    neighbours automatically.
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("{{ c_3 }} V") | Out-Null
$rng.Collapse(0)
$rng.InsertParagraphAfter() | Out-Null

$np = $d.Paragraphs.Item(16).Range
$np.Collapse(1)
$np.InsertAfter("Outlier Min/Max")
$np.Font.Bold = $true
$np.Font.Size = 10

$np.Collapse(0)
$np.InsertAfter([char]9)
$np.Font.Bold = $true
$np.Font.Size = 10

$np.Collapse(0)
$np.InsertAfter([char]9)
$np.Font.Bold = $true
$np.Font.Size = 10

$np.Collapse(0)
$np.InsertAfter("{{ c_4 }}/ {{c_5}}")
$np.Font.Bold = $false
$np.Font.Size = 10

# ---------------------------------------------------------------
# 4. "Total Samples Outside Range:  ~12  " becomes
#    "Total Outlier:  {{ c_6 }}" and the _GoBack bookmark re-appears
#    at the very end of this paragraph.
#    (Work within this one paragraph's Range so the short "~12" needle
#    can't accidentally match text elsewhere in the document.)
# ---------------------------------------------------------------
$para = $d.Paragraphs.Item(17)
$rng = $para.Range
$rng.Find.Execute("Total Samples Outside ", $false, $false, $false, $false, $false, $true, 1, $false, "Total ", 2) | Out-Null

$para = $d.Paragraphs.Item(17)
$rng = $para.Range
$rng.Find.Execute("Range:", $false, $false, $false, $false, $false, $true, 1, $false, "Outlier:", 2) | Out-Null

$para = $d.Paragraphs.Item(17)
$rng = $para.Range
$rng.Find.Execute("~12", $false, $false, $false, $false, $false, $true, 1, $false, "{{ ", 2) | Out-Null
$rng.Font.Bold = $false

$rng.Collapse(0)
$rng.InsertAfter("c_6 }}")
$rng.Font.Bold = $false
$rng.Font.Size = 10
$rng.Font.NameFarEast = "SimSun"

# Drop the trailing (now orphaned) tab run that used to follow "~12".
$para = $d.Paragraphs.Item(17)
$rng = $para.Range
$rng.Collapse(0)
$rng.MoveEnd(1, -1) | Out-Null
$rng.MoveStart(1, -1) | Out-Null
Write-Output "trailing char before delete: [$($rng.Text)]"
if ($rng.Text -eq [string][char]9) {
    $rng.Delete()
}

# Re-add the _GoBack bookmark at the end of this paragraph (right before
# the paragraph mark), matching the diff.
$para = $d.Paragraphs.Item(17)
$rng = $para.Range
$rng.Collapse(0)
$rng.MoveEnd(1, -1) | Out-Null
$d.Bookmarks.Add("_GoBack", $rng) | Out-Null

Write-Output "done"
